# B6-PowerPoint.pptx edit (Thu, Apr 09, 2020 9:05:22 PM)
#
# 1. Re-style the three tables (slides 14-16) from the default
#    "Table_0" style {3465B7D1-8973-47F5-8500-B7799D63012B} to the
#    built-in style {E5703BA1-1461-4C62-9EF8-AA12E118C482}.
# 2. Swap the deck's two DrawingML themes: the slide master currently
#    renders with the "Integral" (Red Violet) palette while the notes
#    master keeps the stock "Office Theme" palette. The edit flips
#    that, so the slides/master now use the plain "Office Theme"
#    colours (the notes master keeps whichever palette this host
#    exposes for editing).

$p = $ppt.ActivePresentation

# --- 1. Table styles ------------------------------------------------
$newTableStyle = "{E5703BA1-1461-4C62-9EF8-AA12E118C482}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Theme colour swap -------------------------------------------
# msoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# Target palette = the stock Office "Office Theme" colours (the colours
# that were previously only on the notes-master-linked theme part).
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}

Write-Output "Applied table style + theme colour updates"
